{"js": "// Wrap the two log-message placeholders (\"EventName\" and \"EventStart\") in the\n// \"create event\" log-format paragraph with OLE_LINK bookmark pairs, exactly as\n// Word leaves behind after a copy/paste round-trip of the selected text.\nconst body = context.document.body;\n\n// Locate the unique occurrences of the two placeholder words in the document.\nconst nameResults = body.search(\"EventName\", { matchCase: true, matchWholeWord: false });\nnameResults.load(\"items\");\nconst startResults = body.search(\"EventStart\", { matchCase: true, matchWholeWord: false });\nstartResults.load(\"items\");\nawait context.sync();\n\nif (nameResults.items.length !== 1) {\n  throw new Error(\"Expected exactly one match for 'EventName', found \" + nameResults.items.length);\n}\nif (startResults.items.length !== 1) {\n  throw new Error(\"Expected exactly one match for 'EventStart', found \" + startResults.items.length);\n}\n\n// Bookmark the \"EventName\" placeholder with OLE_LINK1 / OLE_LINK2.\nconst nameRange = nameResults.items[0];\nnameRange.insertBookmark(\"OLE_LINK1\");\nnameRange.insertBookmark(\"OLE_LINK2\");\n\n// Bookmark the \"EventStart\" placeholder with OLE_LINK3 / OLE_LINK4.\nconst startRange = startResults.items[0];\nstartRange.insertBookmark(\"OLE_LINK3\");\nstartRange.insertBookmark(\"OLE_LINK4\");\n\nawait context.sync();\n", "ps1": "# Wrap the two log-message placeholders (\"EventName\" and \"EventStart\") in the\n# \"create event\" log-format paragraph with OLE_LINK bookmark pairs, exactly as\n# Word leaves behind after a copy/paste round-trip of the selected text.\n$d = $word.ActiveDocument\n\n# Bookmark the \"EventName\" placeholder with OLE_LINK1 / OLE_LINK2.\n$rngName1 = $d.Content\n$rngName1.Find.Execute(\"EventName\") | Out-Null\n$d.Bookmarks.Add(\"OLE_LINK1\", $rngName1) | Out-Null\n\n$rngName2 = $d.Content\n$rngName2.Find.Execute(\"EventName\") | Out-Null\n$d.Bookmarks.Add(\"OLE_LINK2\", $rngName2) | Out-Null\n\n# Bookmark the \"EventStart\" placeholder with OLE_LINK3 / OLE_LINK4.\n$rngStart1 = $d.Content\n$rngStart1.Find.Execute(\"EventStart\") | Out-Null\n$d.Bookmarks.Add(\"OLE_LINK3\", $rngStart1) | Out-Null\n\n$rngStart2 = $d.Content\n$rngStart2.Find.Execute(\"EventStart\") | Out-Null\n$d.Bookmarks.Add(\"OLE_LINK4\", $rngStart2) | Out-Null\n"}
